# Add season record columns (Wins, Losses, Ties) to the COL_2011 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new column labels, matching style of existing header cells.
$ws.Cells.Item(1, 30).Value = "Wins"    # AD1
$ws.Cells.Item(1, 31).Value = "Losses"  # AE1
$ws.Cells.Item(1, 32).Value = "Ties"    # AF1

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows (2-57) - season record values for every player.
for ($row = 2; $row -le 57; $row++) {
    $ws.Cells.Item($row, 30).Value = 73  # AD - Wins
    $ws.Cells.Item($row, 31).Value = 89  # AE - Losses
    $ws.Cells.Item($row, 32).Value = 0   # AF - Ties
}
